$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# ---------------------------------------------------------------------------
# 1) New rows 6 & 7 - bring in the same border/format as row 4 (plain bordered)
#    before writing values, so the new cells pick up the existing "thin border"
#    cell style instead of the engine minting brand-new border objects.
# ---------------------------------------------------------------------------
$ws.Range("A4:G4").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A4:G4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Values for the two new test-case rows
# ---------------------------------------------------------------------------
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "Get Login Detail"
$ws.Cells.Item(6,3).Value = "POST"
$ws.Cells.Item(6,4).Value = "{}"
$ws.Cells.Item(6,5).Value = 200
$ws.Cells.Item(6,6).Value = "success;data.email;data.ClarityID;data.WorldAreaID;data.city;data.country;data.state;data.zipCode;data.latestAppVersion;data.OFSC_UserID;data.Name;data.Currency"
$ws.Cells.Item(6,7).Value = "true;Aashish.Kumar@Emerson.com;AD_JT_01;12061;Fridabad;IN;Haryana;121007;4.7.1;348;Kumar, Aashish;1"

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Validate user detail"
$ws.Cells.Item(7,3).Value = "GET"
$ws.Cells.Item(7,4).Value = "{}"
$ws.Cells.Item(7,5).Value = 200
$ws.Cells.Item(7,6).Value = "success;data.alreadyLoggedIn;data.latestAppVersion"
$ws.Cells.Item(7,7).Value = "true;false;4.7.1"

# ---------------------------------------------------------------------------
# 3) Re-format column F/G of rows 2 & 5, plus the new rows 6 & 7, to the
#    wrap-text styles used by the revised test-data sheet.
# ---------------------------------------------------------------------------

# Row 2: F2 loses the centered alignment, keeps top+wrap only; G2 (hyperlink)
# gains top+wrap.
$f2 = $ws.Range("F2")
$f2.HorizontalAlignment = 1        # xlHAlignGeneral - clear the old "center"
$f2.VerticalAlignment = -4160      # xlTop
$f2.WrapText = $true

$g2 = $ws.Range("G2")
$g2.VerticalAlignment = -4160      # xlTop
$g2.WrapText = $true

# Rows 5, 6, 7: F/G become left+top+wrap (same target style). Set it once on
# F5, then fan the finished format out with a formats-only paste so the
# engine dedupes every cell onto the single resulting shared style instead of
# minting a throwaway style per intermediate property assignment.
$f5 = $ws.Range("F5")
$f5.HorizontalAlignment = -4131  # xlLeft
$f5.VerticalAlignment = -4160    # xlTop
$f5.WrapText = $true

$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 135
$ws.Rows.Item(6).RowHeight = 45

# ---------------------------------------------------------------------------
# 5) Column widths for F & G (bestFit -> fixed custom widths)
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 66.30729166666667
$ws.Columns.Item(7).ColumnWidth = 46.022135416666664

# ---------------------------------------------------------------------------
# 6) Selection / active cell moves to G8 (one past the new last row)
# ---------------------------------------------------------------------------
$ws.Range("G8").Select()
